# Auto-generated edit script for Durandal_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific rows
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 1039.875
$ws.Range("I59").Value = 300
$ws.Range("J59").Value = 1145.5714
$ws.Range("K59").Value = 900
$ws.Range("L59").Value = 3436.7142
$ws.Range("M59").Value = -343
$ws.Range("N59").Value = -4550.7142

$ws.Range("H86").Value = 1657.3939
$ws.Range("I86").Value = 1526.8667
$ws.Range("J86").Value = 1766.1666
$ws.Range("K86").Value = 1526.8667
$ws.Range("L86").Value = 1766.1666
$ws.Range("M86").Value = -403.8667
$ws.Range("N86").Value = -4012.1666

$ws.Range("H89").Value = 1657.3939
$ws.Range("I89").Value = 1526.8667
$ws.Range("J89").Value = 1766.1666
$ws.Range("K89").Value = 7634.333500000001
$ws.Range("L89").Value = 8830.833000000001
$ws.Range("M89").Value = -2018.333500000001
$ws.Range("N89").Value = -20062.833

$ws.Range("H116").Value = 9138.666999999999
$ws.Range("I116").Value = 10062.308
$ws.Range("J116").Value = 3135
$ws.Range("K116").Value = 10062.308
$ws.Range("L116").Value = 3135
$ws.Range("M116").Value = -6620.308000000001
$ws.Range("N116").Value = -10019

$ws.Range("H129").Value = 1158
$ws.Range("I129").Value = 314.66666
$ws.Range("J129").Value = 1222.8718
$ws.Range("K129").Value = 943.9999799999999
$ws.Range("L129").Value = 3668.6154
$ws.Range("M129").Value = 4056.00002
$ws.Range("N129").Value = -13668.6154

$ws.Range("H137").Value = 1231.2632
$ws.Range("I137").Value = 846.0357
$ws.Range("K137").Value = 2538.1071
$ws.Range("M137").Value = 11.89289999999983

$ws.Range("H138").Value = 3400.0144
$ws.Range("I138").Value = 2292.8484
$ws.Range("J138").Value = 4414.9165
$ws.Range("K138").Value = 6878.5452
$ws.Range("L138").Value = 13244.7495
$ws.Range("M138").Value = -1738.5452
$ws.Range("N138").Value = -23524.7495

$ws.Range("H140").Value = 97300
$ws.Range("J140").Value = 97300
$ws.Range("L140").Value = 97300
$ws.Range("N140").Value = -107660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 688908.4399999999
$ws.Range("I32").Value = 8947.102000000001
$ws.Range("J32").Value = 8712452
$ws.Range("K32").Value = 8947.102000000001
$ws.Range("L32").Value = 8712452
$ws.Range("M32").Value = -8660.102000000001
$ws.Range("N32").Value = -8713026

$ws.Range("H37").Value = 5587.143
$ws.Range("J37").Value = 8019
$ws.Range("L37").Value = 8019
$ws.Range("N37").Value = -8565

$ws.Range("H61").Value = 2545.6155
$ws.Range("I61").Value = 2503
$ws.Range("J61").Value = 3057
$ws.Range("K61").Value = 2503
$ws.Range("L61").Value = 3057
$ws.Range("M61").Value = -2291
$ws.Range("N61").Value = -3481

$ws.Range("H74").Value = 1084.3
$ws.Range("I74").Value = 1092.1111
$ws.Range("K74").Value = 1092.1111
$ws.Range("M74").Value = -218.1111000000001

$ws.Range("H77").Value = 1084.3
$ws.Range("I77").Value = 1092.1111
$ws.Range("K77").Value = 5460.5555
$ws.Range("M77").Value = -1092.5555

$ws.Range("H122").Value = 3541.4285
$ws.Range("I122").Value = 3432.4358
$ws.Range("J122").Value = 4958.3335
$ws.Range("K122").Value = 10297.3074
$ws.Range("L122").Value = 14875.0005
$ws.Range("M122").Value = -7847.307400000002
$ws.Range("N122").Value = -19775.0005

$ws.Range("H132").Value = 2296.0925
$ws.Range("I132").Value = 1942.1351
$ws.Range("J132").Value = 3066.4707
$ws.Range("K132").Value = 5826.4053
$ws.Range("L132").Value = 9199.4121
$ws.Range("M132").Value = -3296.4053
$ws.Range("N132").Value = -14259.4121

$ws.Range("H136").Value = 2545.6155
$ws.Range("I136").Value = 2503
$ws.Range("J136").Value = 3057
$ws.Range("K136").Value = 7509
$ws.Range("L136").Value = 9171
$ws.Range("M136").Value = -4959
$ws.Range("N136").Value = -14271

$ws.Range("H141").Value = 65772.89999999999
$ws.Range("J141").Value = 65772.89999999999
$ws.Range("L141").Value = 65772.89999999999
$ws.Range("N141").Value = -76132.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 572.73914
$ws.Range("I80").Value = 759.25
$ws.Range("J80").Value = 473.26666
$ws.Range("K80").Value = 759.25
$ws.Range("L80").Value = 473.26666
$ws.Range("M80").Value = 238.75
$ws.Range("N80").Value = -2469.26666

$ws.Range("H83").Value = 572.73914
$ws.Range("I83").Value = 759.25
$ws.Range("J83").Value = 473.26666
$ws.Range("K83").Value = 3796.25
$ws.Range("L83").Value = 2366.3333
$ws.Range("M83").Value = 1195.75
$ws.Range("N83").Value = -12350.3333

$ws.Range("H134").Value = 1103.6666
$ws.Range("I134").Value = 955.75
$ws.Range("J134").Value = 1399.5
$ws.Range("K134").Value = 2867.25
$ws.Range("L134").Value = 4198.5
$ws.Range("M134").Value = -332.25
$ws.Range("N134").Value = -9268.5

$ws.Range("H135").Value = 45862.223
$ws.Range("J135").Value = 45862.223
$ws.Range("L135").Value = 45862.223
$ws.Range("N135").Value = -56002.223

$ws.Range("H138").Value = 50740
$ws.Range("J138").Value = 50740
$ws.Range("L138").Value = 50740
$ws.Range("N138").Value = -61020

$ws.Range("H140").Value = 89800
$ws.Range("J140").Value = 89800
$ws.Range("L140").Value = 89800
$ws.Range("N140").Value = -100160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9224
$ws.Range("J51").Value = 9224
$ws.Range("L51").Value = 9224
$ws.Range("N51").Value = -10696

$ws.Range("H60").Value = 6776.125
$ws.Range("J60").Value = 8054
$ws.Range("L60").Value = 8054
$ws.Range("N60").Value = -9076

$ws.Range("H61").Value = 9224
$ws.Range("J61").Value = 9224
$ws.Range("L61").Value = 9224
$ws.Range("N61").Value = -9920

$ws.Range("H74").Value = 17391.143
$ws.Range("J74").Value = 17391.143
$ws.Range("L74").Value = 17391.143
$ws.Range("N74").Value = -19139.143

$ws.Range("H77").Value = 17391.143
$ws.Range("J77").Value = 17391.143
$ws.Range("L77").Value = 52173.429
$ws.Range("N77").Value = -60909.429

$ws.Range("H107").Value = 603.74286
$ws.Range("I107").Value = 453.54544
$ws.Range("K107").Value = 453.54544
$ws.Range("M107").Value = 1466.45456

$ws.Range("H115").Value = 49980
$ws.Range("J115").Value = 49980
$ws.Range("L115").Value = 49980
$ws.Range("N115").Value = -52330

$ws.Range("H122").Value = 11766109
$ws.Range("I122").Value = 1249
$ws.Range("J122").Value = 22223762
$ws.Range("K122").Value = 3747
$ws.Range("L122").Value = 66671286
$ws.Range("M122").Value = -1297
$ws.Range("N122").Value = -66676186

$ws.Range("H141").Value = 34500
$ws.Range("I141").Value = 19000
$ws.Range("K141").Value = 19000
$ws.Range("M141").Value = -13820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 217.90909
$ws.Range("J12").Value = 217.90909
$ws.Range("L12").Value = 653.72727
$ws.Range("N12").Value = -999.72727

$ws.Range("H114").Value = 306785.7
$ws.Range("I114").Value = 7621.8667
$ws.Range("J114").Value = 556088.9
$ws.Range("K114").Value = 22865.6001
$ws.Range("L114").Value = 1668266.7
$ws.Range("M114").Value = -19611.6001
$ws.Range("N114").Value = -1674774.7

$ws.Range("H117").Value = 85133.336
$ws.Range("J117").Value = 101890
$ws.Range("L117").Value = 305670
$ws.Range("N117").Value = -312554

$ws.Range("H121").Value = 51409.2
$ws.Range("J121").Value = 67451.07000000001
$ws.Range("L121").Value = 202353.21
$ws.Range("N121").Value = -204973.21

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2682.8572
$ws.Range("I80").Value = 2630
$ws.Range("K80").Value = 2630
$ws.Range("M80").Value = -1632

$ws.Range("H83").Value = 2682.8572
$ws.Range("I83").Value = 2630
$ws.Range("K83").Value = 13150
$ws.Range("M83").Value = -8158

$ws.Range("H102").Value = 1752.742
$ws.Range("I102").Value = 1470.96
$ws.Range("K102").Value = 1470.96
$ws.Range("M102").Value = 151.04

$ws.Range("H113").Value = 4584984
$ws.Range("I113").Value = 8334596
$ws.Range("K113").Value = 8334596
$ws.Range("M113").Value = -8332426

$ws.Range("H131").Value = 35326
$ws.Range("J131").Value = 35326
$ws.Range("L131").Value = 35326
$ws.Range("N131").Value = -45406

$ws.Range("H138").Value = 69600
$ws.Range("J138").Value = 69600
$ws.Range("L138").Value = 69600
$ws.Range("N138").Value = -79880

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3981.25
$ws.Range("I7").Value = 3334.625
$ws.Range("J7").Value = 4627.875
$ws.Range("K7").Value = 3334.625
$ws.Range("L7").Value = 4627.875
$ws.Range("M7").Value = -3222.625
$ws.Range("N7").Value = -4851.875

$ws.Range("H126").Value = 3981.25
$ws.Range("I126").Value = 3334.625
$ws.Range("J126").Value = 4627.875
$ws.Range("K126").Value = 10003.875
$ws.Range("L126").Value = 13883.625
$ws.Range("M126").Value = -7533.875
$ws.Range("N126").Value = -18823.625

$ws.Range("H136").Value = 3352.1943
$ws.Range("I136").Value = 1921.76
$ws.Range("J136").Value = 6603.1816
$ws.Range("K136").Value = 5765.28
$ws.Range("L136").Value = 19809.5448
$ws.Range("M136").Value = -3215.28
$ws.Range("N136").Value = -24909.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2517.3333
$ws.Range("I107").Value = 2660.8
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 7982.400000000001
$ws.Range("L107").Value = 5400
$ws.Range("M107").Value = -6062.400000000001
$ws.Range("N107").Value = -9240

$ws.Range("H138").Value = 69533.336
$ws.Range("J138").Value = 86800
$ws.Range("L138").Value = 86800
$ws.Range("N138").Value = -97080

Write-Output "Updated 47 rows across 8 sheets"